# Insert 3 new rows before row 331 (a new weekly price group for
# "Comercializadora del Agro de Limarí" / Pepino dulce, date 44615,
# qualities Primera/Segunda/Tercera). This pushes the existing rows
# 331-344 down to 334-347, matching the target dimension A1:R347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A331:A333").EntireRow.Insert()

# Row 331: Primera
$ws.Cells.Item(331,1).Value = 2
$ws.Cells.Item(331,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(331,3).Value = "Coquimbo"
$ws.Cells.Item(331,4).Value = 44615
$ws.Cells.Item(331,5).Value = 4
$ws.Cells.Item(331,6).Value = 100112043
$ws.Cells.Item(331,7).Value = "Pepino dulce"
$ws.Cells.Item(331,8).Value = "Cultivar IV Región"
$ws.Cells.Item(331,9).Value = "Primera"
$ws.Cells.Item(331,10).Value = 400
$ws.Cells.Item(331,11).Value = 10500
$ws.Cells.Item(331,12).Value = 11000
$ws.Cells.Item(331,13).Value = 10750
$ws.Cells.Item(331,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(331,15).Value = "Provincia de Limarí"
$ws.Cells.Item(331,16).Value = 597
$ws.Cells.Item(331,17).Value = 18
$ws.Cells.Item(331,18).Value = "Hortaliza"

# Row 332: Segunda
$ws.Cells.Item(332,1).Value = 2
$ws.Cells.Item(332,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(332,3).Value = "Coquimbo"
$ws.Cells.Item(332,4).Value = 44615
$ws.Cells.Item(332,5).Value = 4
$ws.Cells.Item(332,6).Value = 100112043
$ws.Cells.Item(332,7).Value = "Pepino dulce"
$ws.Cells.Item(332,8).Value = "Cultivar IV Región"
$ws.Cells.Item(332,9).Value = "Segunda"
$ws.Cells.Item(332,10).Value = 400
$ws.Cells.Item(332,11).Value = 8500
$ws.Cells.Item(332,12).Value = 9000
$ws.Cells.Item(332,13).Value = 8750
$ws.Cells.Item(332,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(332,15).Value = "Provincia de Limarí"
$ws.Cells.Item(332,16).Value = 486
$ws.Cells.Item(332,17).Value = 18
$ws.Cells.Item(332,18).Value = "Hortaliza"

# Row 333: Tercera
$ws.Cells.Item(333,1).Value = 2
$ws.Cells.Item(333,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(333,3).Value = "Coquimbo"
$ws.Cells.Item(333,4).Value = 44615
$ws.Cells.Item(333,5).Value = 4
$ws.Cells.Item(333,6).Value = 100112043
$ws.Cells.Item(333,7).Value = "Pepino dulce"
$ws.Cells.Item(333,8).Value = "Cultivar IV Región"
$ws.Cells.Item(333,9).Value = "Tercera"
$ws.Cells.Item(333,10).Value = 300
$ws.Cells.Item(333,11).Value = 6500
$ws.Cells.Item(333,12).Value = 7000
$ws.Cells.Item(333,13).Value = 6750
$ws.Cells.Item(333,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(333,15).Value = "Provincia de Limarí"
$ws.Cells.Item(333,16).Value = 375
$ws.Cells.Item(333,17).Value = 18
$ws.Cells.Item(333,18).Value = "Hortaliza"
